# Insert a new weekly data row for "Cebollín" at row 155, shifting the
# existing rows 155:182 down to 156:183 (dimension grows to A1:R183).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(155).Insert()

$ws.Range("A155").Value = 7
$ws.Range("B155").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C155").Value = "Ñuble"
$ws.Range("D155").Value = 45131
$ws.Range("E155").Value = 16
$ws.Range("F155").Value = 100112037
$ws.Range("G155").Value = "Cebollín"
$ws.Range("H155").Value = "Sin especificar"
$ws.Range("I155").Value = "Primera"
$ws.Range("J155").Value = 100
$ws.Range("K155").Value = 7000
$ws.Range("L155").Value = 7000
$ws.Range("M155").Value = 7000
$ws.Range("N155").Value = "`$/paquete 36 unidades"
$ws.Range("O155").Value = "Provincia de Diguillín"
$ws.Range("P155").Value = 194
$ws.Range("Q155").Value = 36
$ws.Range("R155").Value = "Hortaliza"
